$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update header row text (row 1) to the new test-case-sheet terminology
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Test_case_ID"
$ws.Range("B1").Value = "Test Scenario"
$ws.Range("C1").Value = "Pre-Condition"
$ws.Range("D1").Value = "Test Steps"
# E1 "Test Data" is unchanged
$ws.Range("F1").Value = "Expected Result"
$ws.Range("G1").Value = "Severity"
# H1 "Priority" is unchanged

# ---------------------------------------------------------------------------
# 2) Re-shuffle cell formatting (direct styles) between header/row2 cells.
#    Capture the "donor" styles that are about to move before overwriting
#    the cells that currently hold them.
# ---------------------------------------------------------------------------

# C1 currently carries the "black font" style -> needed on C2 afterwards.
$ws.Range("C1").Copy()
$ws.Range("C2").PasteSpecial(-4122)

# E1/G1/H1 currently carry the "dark-grey font" style -> needed on E2/G2/H2.
$ws.Range("E1").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("E1").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("E1").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# D1 keeps the plain/default style throughout - use it as the donor for all
# cells that need to fall back to the default style.
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("D1").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# B3 keeps the "wrap text" style throughout - use it as the donor for the
# cells that need to pick that style up.
$ws.Range("B3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Sheet view: drop the frozen/scrolled "topLeftCell=E1" viewport and
#    narrow the active selection down from the whole column to just H1.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("H1").Select() | Out-Null
